# Auto-generated edit script: refresh market-price derived values
# in the Seraph_Profits workbook (Leve profit calcs), per sheet/row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 583.6061
$ws.Range("I15").Value2 = 583.6061
$ws.Range("K15").Value2 = 1750.8183
$ws.Range("M15").Value2 = -1581.8183
$ws.Range("H18").Value2 = 1909.9
$ws.Range("I18").Value2 = 1987.375
$ws.Range("K18").Value2 = 1987.375
$ws.Range("M18").Value2 = -1703.375
$ws.Range("H28").Value2 = 514.5
$ws.Range("I28").Value2 = 620.7143
$ws.Range("K28").Value2 = 620.7143
$ws.Range("M28").Value2 = -135.7143
$ws.Range("H40").Value2 = 2381.4546
$ws.Range("J40").Value2 = 3250
$ws.Range("L40").Value2 = 3250
$ws.Range("N40").Value2 = -3600
$ws.Range("H43").Value2 = 3982.3333
$ws.Range("I43").Value2 = 1800
$ws.Range("J43").Value2 = 6164.6665
$ws.Range("K43").Value2 = 1800
$ws.Range("L43").Value2 = 6164.6665
$ws.Range("M43").Value2 = -1731
$ws.Range("N43").Value2 = -6302.6665
$ws.Range("H107").Value2 = 497
$ws.Range("I107").Value2 = 496.7
$ws.Range("K107").Value2 = 496.7
$ws.Range("M107").Value2 = 1423.3
$ws.Range("H138").Value2 = 4546.905
$ws.Range("I138").Value2 = 1995.5714
$ws.Range("K138").Value2 = 5986.7142
$ws.Range("M138").Value2 = -846.7142000000003
$ws.Range("H141").Value2 = 6331.6665
$ws.Range("I141").Value2 = 5997.5
$ws.Range("K141").Value2 = 17992.5
$ws.Range("M141").Value2 = -12812.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 16650.797
$ws.Range("I32").Value2 = 7450.4194
$ws.Range("K32").Value2 = 7450.4194
$ws.Range("M32").Value2 = -7163.4194
$ws.Range("H63").Value2 = 5146.722
$ws.Range("I63").Value2 = 4037.125
$ws.Range("J63").Value2 = 6034.4
$ws.Range("K63").Value2 = 4037.125
$ws.Range("L63").Value2 = 6034.4
$ws.Range("M63").Value2 = -3351.125
$ws.Range("N63").Value2 = -7406.4
$ws.Range("H66").Value2 = 5146.722
$ws.Range("I66").Value2 = 4037.125
$ws.Range("J66").Value2 = 6034.4
$ws.Range("K66").Value2 = 20185.625
$ws.Range("L66").Value2 = 30172
$ws.Range("M66").Value2 = -16753.625
$ws.Range("N66").Value2 = -37036
$ws.Range("H74").Value2 = 8840.333000000001
$ws.Range("I74").Value2 = 0
$ws.Range("K74").Value2 = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value2 = 8840.333000000001
$ws.Range("I77").Value2 = 0
$ws.Range("K77").Value2 = 0
$ws.Range("M77").ClearContents()
$ws.Range("H132").Value2 = 2083.5
$ws.Range("J132").Value2 = 2996
$ws.Range("L132").Value2 = 8988
$ws.Range("N132").Value2 = -14048

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 1531.8334
$ws.Range("I20").Value2 = 1470.4286
$ws.Range("J20").Value2 = 1617.8
$ws.Range("K20").Value2 = 1470.4286
$ws.Range("L20").Value2 = 1617.8
$ws.Range("M20").Value2 = -1223.4286
$ws.Range("N20").Value2 = -2111.8
$ws.Range("H99").Value2 = 974.5
$ws.Range("I99").Value2 = 974.5
$ws.Range("K99").Value2 = 974.5
$ws.Range("M99").Value2 = 523.5
$ws.Range("H107").Value2 = 934.2857
$ws.Range("I107").Value2 = 934.2857
$ws.Range("K107").Value2 = 934.2857
$ws.Range("M107").Value2 = 985.7143
$ws.Range("H134").Value2 = 4252.75
$ws.Range("I134").Value2 = 3870.3333
$ws.Range("K134").Value2 = 11610.9999
$ws.Range("M134").Value2 = -9075.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 4769.3
$ws.Range("I31").Value2 = 3224.4
$ws.Range("J31").Value2 = 6314.2
$ws.Range("K31").Value2 = 3224.4
$ws.Range("L31").Value2 = 6314.2
$ws.Range("M31").Value2 = -2929.4
$ws.Range("N31").Value2 = -6904.2
$ws.Range("H34").Value2 = 4769.3
$ws.Range("I34").Value2 = 3224.4
$ws.Range("J34").Value2 = 6314.2
$ws.Range("K34").Value2 = 3224.4
$ws.Range("L34").Value2 = 6314.2
$ws.Range("M34").Value2 = -3022.4
$ws.Range("N34").Value2 = -6718.2
$ws.Range("H107").Value2 = 312.2143
$ws.Range("I107").Value2 = 240.07692
$ws.Range("J107").Value2 = 1250
$ws.Range("K107").Value2 = 240.07692
$ws.Range("L107").Value2 = 1250
$ws.Range("M107").Value2 = 1679.92308
$ws.Range("N107").Value2 = -5090
$ws.Range("H132").Value2 = 4229.364
$ws.Range("J132").Value2 = 5404.6
$ws.Range("L132").Value2 = 16213.8
$ws.Range("N132").Value2 = -21273.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 105964024
$ws.Range("I4").Value2 = 111850420
$ws.Range("J4").Value2 = 8888
$ws.Range("K4").Value2 = 335551260
$ws.Range("L4").Value2 = 26664
$ws.Range("M4").Value2 = -335551148
$ws.Range("N4").Value2 = -26888
$ws.Range("H14").Value2 = 1059.9166
$ws.Range("I14").Value2 = 1059.9166
$ws.Range("K14").Value2 = 3179.7498
$ws.Range("M14").Value2 = -3006.7498
$ws.Range("H21").Value2 = 295
$ws.Range("I21").Value2 = 0
$ws.Range("J21").Value2 = 295
$ws.Range("K21").Value2 = 0
$ws.Range("L21").Value2 = 885
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value2 = -1231
$ws.Range("H23").Value2 = 151.83333
$ws.Range("I23").Value2 = 232.33333
$ws.Range("K23").Value2 = 696.99999
$ws.Range("M23").Value2 = -461.99999
$ws.Range("H75").Value2 = 372.14285
$ws.Range("I75").Value2 = 392.5
$ws.Range("K75").Value2 = 1177.5
$ws.Range("M75").Value2 = -179.5
$ws.Range("H78").Value2 = 372.14285
$ws.Range("I78").Value2 = 392.5
$ws.Range("K78").Value2 = 3532.5
$ws.Range("M78").Value2 = 1459.5
$ws.Range("H81").Value2 = 1762.4
$ws.Range("J81").Value2 = 1725.5
$ws.Range("L81").Value2 = 5176.5
$ws.Range("N81").Value2 = -7422.5
$ws.Range("H84").Value2 = 1762.4
$ws.Range("J84").Value2 = 1725.5
$ws.Range("L84").Value2 = 15529.5
$ws.Range("N84").Value2 = -26761.5
$ws.Range("H92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("L92").Value2 = 0
$ws.Range("N92").ClearContents()
$ws.Range("H109").Value2 = 1296.5
$ws.Range("I109").Value2 = 445
$ws.Range("K109").Value2 = 1335
$ws.Range("M109").Value2 = -295
$ws.Range("H121").Value2 = 908.5714
$ws.Range("I121").Value2 = 226.66667
$ws.Range("K121").Value2 = 680.00001
$ws.Range("M121").Value2 = 629.99999
$ws.Range("H129").Value2 = 2618.6155
$ws.Range("I129").Value2 = 1629.4
$ws.Range("J129").Value2 = 3236.875
$ws.Range("K129").Value2 = 4888.200000000001
$ws.Range("L129").Value2 = 9710.625
$ws.Range("M129").Value2 = 111.7999999999993
$ws.Range("N129").Value2 = -19710.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 7197
$ws.Range("I132").Value2 = 6633.3335
$ws.Range("K132").Value2 = 19900.0005
$ws.Range("M132").Value2 = -17370.0005
$ws.Range("H137").Value2 = 30354.5
$ws.Range("I137").Value2 = 30354.5
$ws.Range("K137").Value2 = 30354.5
$ws.Range("M137").Value2 = -25254.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 493.33334
$ws.Range("I22").Value2 = 523.75
$ws.Range("K22").Value2 = 523.75
$ws.Range("M22").Value2 = -228.75
$ws.Range("H27").Value2 = 493.33334
$ws.Range("I27").Value2 = 523.75
$ws.Range("K27").Value2 = 523.75
$ws.Range("M27").Value2 = -416.75
$ws.Range("H61").Value2 = 2692.5557
$ws.Range("I61").Value2 = 2364.2666
$ws.Range("K61").Value2 = 2364.2666
$ws.Range("M61").Value2 = -2162.2666
$ws.Range("H113").Value2 = 2692.5557
$ws.Range("I113").Value2 = 2364.2666
$ws.Range("K113").Value2 = 2364.2666
$ws.Range("M113").Value2 = -194.2665999999999
$ws.Range("H132").Value2 = 5395.241
$ws.Range("I132").Value2 = 4959.3335
$ws.Range("J132").Value2 = 6108.5454
$ws.Range("K132").Value2 = 14878.0005
$ws.Range("L132").Value2 = 18325.6362
$ws.Range("M132").Value2 = -12348.0005
$ws.Range("N132").Value2 = -23385.6362
$ws.Range("H136").Value2 = 3500.75
$ws.Range("I136").Value2 = 3167.6667
$ws.Range("J136").Value2 = 4500
$ws.Range("K136").Value2 = 9503.000100000001
$ws.Range("L136").Value2 = 13500
$ws.Range("M136").Value2 = -6953.000100000001
$ws.Range("N136").Value2 = -18600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 1171.091
$ws.Range("I100").Value2 = 1081.8889
$ws.Range("J100").Value2 = 1572.5
$ws.Range("K100").Value2 = 2163.7778
$ws.Range("L100").Value2 = 3145
$ws.Range("M100").Value2 = -1622.7778
$ws.Range("N100").Value2 = -4227
$ws.Range("H132").Value2 = 2004
$ws.Range("I132").Value2 = 2004
$ws.Range("K132").Value2 = 6012
$ws.Range("M132").Value2 = -3482
$ws.Range("H136").Value2 = 54447.473
$ws.Range("I136").Value2 = 1522.7273
$ws.Range("K136").Value2 = 4568.1819
$ws.Range("M136").Value2 = -2018.1819
$ws.Range("H140").Value2 = 99900
$ws.Range("J140").Value2 = 99900
$ws.Range("L140").Value2 = 99900
$ws.Range("N140").Value2 = -110260
